$d = $word.ActiveDocument
$d.Content.Find.Execute("For eature selection for objective one", $true, $false, $false, $false, $false,
                         $true, 1, $false, "For feature selection for objective one", 2)
